$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 11) for the two
# styled columns (A = id-like indexed column, G = date column) so the new
# row matches the workbook's existing per-column styles, then overwrite
# the values that were copied along with the formatting.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("G11").Copy($ws.Range("G12"))

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "yes"
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "yes"
$ws.Range("G12").Value = 43723
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 30
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = "no"
$ws.Range("L12").Value = "meditation"
$ws.Range("M12").Value = 12
